$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Revert job image references in column A back to "blog_1.jpg"
$ws.Range("A2").Value = "blog_1.jpg"
$ws.Range("A4").Value = "blog_1.jpg"

# Restore original selection to A2
$null = $ws.Range("A2").Select()
